$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric values for rows 2-6 (columns D:AJ) ---
# Row 2
$ws.Range("D2").Value = 118
$ws.Range("E2").Value = -61
$ws.Range("F2").Value = -61
$ws.Range("G2").Value = -90
$ws.Range("H2").Value = -90
$ws.Range("I2").Value = -88
$ws.Range("J2").Value = -2
$ws.Range("K2").Value = 197
$ws.Range("L2").Value = 114
$ws.Range("M2").Value = 83
$ws.Range("N2").Value = 82
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 42
$ws.Range("Q2").Value = -52
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 54
$ws.Range("T2").Value = 7
$ws.Range("U2").Value = -60
$ws.Range("V2").Value = 72
$ws.Range("W2").Value = -51.58
$ws.Range("X2").Value = -76.84999999999999
$ws.Range("Y2").Value = -74.13
$ws.Range("Z2").Value = -43.42
$ws.Range("AA2").Value = 136.59
$ws.Range("AB2").Value = 120.87
$ws.Range("AC2").Value = -924
$ws.Range("AD2").Value = -3.61
$ws.Range("AE2").Value = 815
$ws.Range("AF2").Value = 4.09
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 10023862

# Row 3
$ws.Range("D3").Value = 115
$ws.Range("E3").Value = -75
$ws.Range("F3").Value = -75
$ws.Range("G3").Value = -78
$ws.Range("H3").Value = -78
$ws.Range("I3").Value = -76
$ws.Range("J3").Value = -2
$ws.Range("K3").Value = 164
$ws.Range("L3").Value = 102
$ws.Range("M3").Value = 62
$ws.Range("N3").Value = 62
$ws.Range("P3").Value = 54
$ws.Range("Q3").Value = -36
$ws.Range("R3").Value = -11
$ws.Range("S3").Value = 41
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = -50
$ws.Range("V3").Value = 30
$ws.Range("W3").Value = -65.52
$ws.Range("X3").Value = -67.70999999999999
$ws.Range("Y3").Value = -105.83
$ws.Range("Z3").Value = -42.98
$ws.Range("AA3").Value = 165.22
$ws.Range("AB3").Value = 39.82
$ws.Range("AC3").Value = -661
$ws.Range("AD3").Value = -5.45
$ws.Range("AE3").Value = 482
$ws.Range("AF3").Value = 7.48
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 12839197

# Row 4
$ws.Range("D4").Value = 83
$ws.Range("E4").Value = -48
$ws.Range("F4").Value = -57
$ws.Range("G4").Value = -69
$ws.Range("H4").Value = -86
$ws.Range("I4").Value = -86
$ws.Range("K4").Value = 288
$ws.Range("L4").Value = 154
$ws.Range("M4").Value = 134
$ws.Range("N4").Value = 134
$ws.Range("P4").Value = 92
$ws.Range("Q4").Value = -49
$ws.Range("R4").Value = -153
$ws.Range("S4").Value = 202
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = -50
$ws.Range("V4").Value = 87
$ws.Range("W4").Value = -57.53
$ws.Range("X4").Value = -103.71
$ws.Range("Y4").Value = -87.87
$ws.Range("Z4").Value = -38.1
$ws.Range("AA4").Value = 114.68
$ws.Range("AB4").Value = 48.59
$ws.Range("AC4").Value = -507
$ws.Range("AD4").Value = -21.5
$ws.Range("AE4").Value = 729
$ws.Range("AF4").Value = 14.95
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 18429147

# Row 5
$ws.Range("D5").Value = 356
$ws.Range("E5").Value = -31
$ws.Range("F5").Value = -31
$ws.Range("G5").Value = -26
$ws.Range("H5").Value = -29
$ws.Range("I5").Value = -29
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 367
$ws.Range("L5").Value = 180
$ws.Range("M5").Value = 187
$ws.Range("N5").Value = 185
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 108
$ws.Range("Q5").Value = -127
$ws.Range("R5").Value = 132
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 1
$ws.Range("U5").Value = -128
$ws.Range("V5").Value = 88
$ws.Range("W5").Value = -8.75
$ws.Range("X5").Value = -8.16
$ws.Range("Y5").Value = -18.2
$ws.Range("Z5").Value = -8.859999999999999
$ws.Range("AA5").Value = 96.38
$ws.Range("AB5").Value = 67.70999999999999
$ws.Range("AC5").Value = -148
$ws.Range("AD5").Value = -23.81
$ws.Range("AE5").Value = 859
$ws.Range("AF5").Value = 4.1
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 21553594

# Row 6
$ws.Range("D6").Value = 345
$ws.Range("E6").Value = -66
$ws.Range("F6").Value = -66
$ws.Range("G6").Value = -117
$ws.Range("H6").Value = -121
$ws.Range("I6").Value = -120
$ws.Range("K6").Value = 614
$ws.Range("L6").Value = 429
$ws.Range("M6").Value = 185
$ws.Range("N6").Value = 184
$ws.Range("P6").Value = 121
$ws.Range("Q6").Value = -27
$ws.Range("R6").Value = -5
$ws.Range("S6").Value = 370
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = -28
$ws.Range("V6").Value = 304
$ws.Range("W6").Value = -19.23
$ws.Range("X6").Value = -35.13
$ws.Range("Y6").Value = -65.25
$ws.Range("Z6").Value = -24.67
$ws.Range("AA6").Value = 231.91
$ws.Range("AB6").Value = 48.94
$ws.Range("AC6").Value = -531
$ws.Range("AD6").Value = -4.04
$ws.Range("AE6").Value = 757
$ws.Range("AF6").Value = 2.83
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AJ6").Value = 24298019

# --- Clear cells removed from rows 3, 4, 6 ---
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("AI6").ClearContents()

# --- Clear all data (D:AI) for rows 7, 8, 9, keeping A, B, C ---
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
